# Fixed #366 User content is lost after two generation without edition.
#
# The "user content" markers (m:usercontent / m:endusercontent) used to
# be stored as simple fields (<w:fldSimple w:instr="..."/>). A simple
# field has no separate "result" part, so the next M2Doc generation
# could not locate/preserve any content a user had typed between the
# markers. Storing the markers as complex fields instead
# (fldChar begin/instrText/fldChar separate/fldChar end) gives each
# field a real result range that a later generation can read back and
# keep, instead of the field - and any user edits next to it - being
# lost.
#
# This script converts every simple field of the document into the
# equivalent complex field, leaving the instruction text and the
# hosting paragraph otherwise untouched.

$d = $word.ActiveDocument

function Get-ParagraphStartContaining($pos) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $candidate = $d.Paragraphs.Item($i)
        if ($pos -ge $candidate.Range.Start -and $pos -lt $candidate.Range.End) {
            return $candidate.Range.Start
        }
    }
    return $pos
}

$complexFieldTemplate = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:fldChar w:fldCharType="begin"/></w:r>
<w:r><w:instrText>__INSTR__</w:instrText></w:r>
<w:r><w:fldChar w:fldCharType="separate"/></w:r>
<w:r><w:fldChar w:fldCharType="end"/></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

# Snapshot the simple fields that exist right now (their instruction
# text and the offset of the field code) before anything is touched:
# converting a field to its complex form registers two *new* fields
# (the begin/.../end run sequence is itself seen as a field), so the
# collection must not be walked/re-queried while mutating it.
$targets = New-Object System.Collections.ArrayList
foreach ($f in $d.Fields) {
    [void]$targets.Add(@{ Instr = $f.Code.Text; Pos = $f.Code.Start })
}

foreach ($target in $targets) {
    $paraStart = Get-ParagraphStartContaining $target.Pos

    # Find the field currently anchored at that paragraph and drop its
    # <w:fldSimple/> representation, leaving the (still intact)
    # paragraph as the insertion anchor.
    foreach ($f in $d.Fields) {
        if ($f.Code.Start -eq $target.Pos) {
            $f.Delete()
            break
        }
    }

    $xml = $complexFieldTemplate.Replace("__INSTR__", $target.Instr)
    $insertionPoint = $d.Range($paraStart, $paraStart)
    $insertionPoint.InsertXML($xml)
}
